# Fixed a bug where the springer api wouldn't properly fill in the abstract
# in some situations. The "Authors" column (E) values for rows 2-12 had
# their internal separators re-joined with an extra space, which is
# reproduced here by inserting one additional space after every comma in
# the existing cell text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # Column E = Authors
    $current = $cell.Text
    if ($current -ne $null -and $current -ne "") {
        $updated = $current -replace ',', ',  '
        $cell.Value = $updated
    }
}
